$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "M2"  = "[49.81015528245323, 50.08267168910147]"
    "U2"  = "[49.86131755371875, 50.049134271255454]"
    "M3"  = "[49.85183558924748, 50.15563036274177]"
    "U3"  = "[49.95073208067164, 50.11446128106226]"
    "M4"  = "[49.94788814097688, 50.23117095913138]"
    "U4"  = "[49.90100364296392, 50.05272108274928]"
    "M5"  = "[49.84571343736905, 50.11868049600594]"
    "U5"  = "[49.94803758169328, 50.1059246608149]"
    "M6"  = "[49.9471566435129, 50.18397158579927]"
    "U6"  = "[49.89451467942269, 50.056461976637515]"
    "M7"  = "[49.82348890073211, 50.1170333951301]"
    "U7"  = "[49.915138467005455, 50.076184932392174]"
    "M8"  = "[49.82918872758348, 50.152036301685975]"
    "U8"  = "[49.91821791307012, 50.08893680237466]"
    "M9"  = "[49.808589612522795, 50.10038479392148]"
    "U9"  = "[49.91323125039128, 50.07561312992928]"
    "M10" = "[50.004311065489325, 50.23589819803967]"
    "U10" = "[49.90952232740756, 50.0635787929066]"
    "M11" = "[49.72734805167844, 50.071070582671545]"
    "U11" = "[49.8771789927341, 50.0605497256886]"
    "M12" = "[49.85223464601765, 50.19856053787758]"
    "U12" = "[49.93030652603219, 50.10750959686326]"
    "M13" = "[49.94492463672572, 50.19940021723759]"
    "U13" = "[49.89102295926903, 50.044119526022904]"
    "M14" = "[49.87554831092492, 50.13450242309933]"
    "U14" = "[49.97364005219648, 50.135333632559274]"
    "M15" = "[49.86684778423538, 50.18043232663087]"
    "U15" = "[49.901138444052684, 50.06455416760681]"
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
